$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "46.052.61"
Set-TextValue "E2" "  -1.75%  "
Set-TextValue "D3" "2.351.96"
Set-TextValue "E3" "  +1.78%  "
Set-TextValue "D4" "0.998"
Set-TextValue "E4" "  -0.17%  "
Set-TextValue "D5" "301.03"
Set-TextValue "E5" "  +0.63%  "
Set-TextValue "D6" "99.43"
Set-TextValue "E6" "  +1.06%  "
Set-TextValue "D7" "0.570"
Set-TextValue "E7" "  -0.55%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  +0.01%  "
Set-TextValue "E9" "  -3.33%  "
Set-TextValue "D10" "34.35"
Set-TextValue "E10" "  -3.68%  "
Set-TextValue "D11" "0.0796"
Set-TextValue "E11" "  -0.24%  "
Set-TextValue "D12" "7.11"
Set-TextValue "E12" "  -3.14%  "
Set-TextValue "E13" "  -0.38%  "
Set-TextValue "D14" "2.709.55"
Set-TextValue "E14" "  +1.77%  "
Set-TextValue "D15" "2.353.82"
Set-TextValue "E15" "  +1.77%  "
Set-TextValue "B16" "Chainlink"
Set-TextValue "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D16" "13.63"
Set-TextValue "E16" "  -2.39%  "
Set-TextValue "B17" "Polygon"
Set-TextValue "C17" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D17" "0.808"
Set-TextValue "E17" "  -1.51%  "
Set-TextValue "D18" "45.929.86"
Set-TextValue "E18" "  -1.68%  "
Set-TextValue "D19" "12.75"
Set-TextValue "E19" "  -3.65%  "
Set-TextValue "D20" "0.0₃0965"
Set-TextValue "E20" "  +2.51%  "
Set-TextValue "D21" "6.01"
Set-TextValue "E21" "  -1.94%  "
Set-TextValue "D22" "67.32"
Set-TextValue "E22" "  +0.59%  "
Set-TextValue "D23" "245.94"
Set-TextValue "E23" "  -1.35%  "
Set-TextValue "D24" "2.83"
Set-TextValue "E24" "  -3.21%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  +0.15%  "
Set-TextValue "D26" "1.92"
Set-TextValue "E26" "  -3.28%  "
Set-TextValue "D27" "39.81"
Set-TextValue "D28" "2.19"
Set-TextValue "E28" "  -2.93%  "
Set-TextValue "D29" "9.78"
Set-TextValue "D30" "20.92"
Set-TextValue "E30" "  +3.68%  "
Set-TextValue "D31" "3.69"
Set-TextValue "E31" "  +18.10%  "
Set-TextValue "E32" "  +5.49%  "
Set-TextValue "E33" "  -4.04%  "
Set-TextValue "D34" "146.18"
Set-TextValue "E34" "  -0.55%  "
Set-TextValue "E35" "  -3.12%  "
Set-TextValue "E36" "  +0.17%  "
Set-TextValue "D37" "1.87"
Set-TextValue "E37" "  +4.04%  "
Set-TextValue "E38" "  -2.10%  "
Set-TextValue "D39" "15.00"
Set-TextValue "E39" "  -3.83%  "
Set-TextValue "D40" "3.94"
Set-TextValue "E40" "  -1.63%  "
Set-TextValue "E41" "  -2.62%  "
Set-TextValue "E42" "  -7.12%  "
Set-TextValue "D43" "1.876.41"
Set-TextValue "E43" "  +2.26%  "
Set-TextValue "D44" "0.998"
Set-TextValue "D45" "93.03"
Set-TextValue "E45" "  +2.03%  "
Set-TextValue "E46" "  -10.97%  "
Set-TextValue "E47" "  -6.67%  "
Set-TextValue "D48" "8.25"
Set-TextValue "E48" "  +3.56%  "
Set-TextValue "D49" "97.59"
Set-TextValue "E49" "  +0.24%  "
Set-TextValue "D50" "2.580.54"
Set-TextValue "E50" "  +1.53%  "
Set-TextValue "B51" "ordi"
Set-TextValue "C51" "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue "D51" "68.48"
Set-TextValue "E51" "  -9.67%  "
